$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 607.1429000000001
$ws.Range("I2").Value = 125
$ws.Range("J2").Value = 1250
$ws.Range("K2").Value = 125
$ws.Range("L2").Value = 1250
$ws.Range("M2").Value = -12
$ws.Range("N2").Value = -1476
$ws.Range("H4").Value = 76.333336
$ws.Range("I4").Value = 76.333336
$ws.Range("K4").Value = 76.333336
$ws.Range("M4").Value = 37.666664
$ws.Range("H17").Value = 4264.1113
$ws.Range("I17").Value = 437.5
$ws.Range("K17").Value = 1312.5
$ws.Range("M17").Value = -1144.5
$ws.Range("H33").Value = 180.57143
$ws.Range("I33").Value = 180.57143
$ws.Range("K33").Value = 180.57143
$ws.Range("M33").Value = 48.42857000000001
$ws.Range("H62").Value = 6151.375
$ws.Range("I62").Value = 4733.1665
$ws.Range("K62").Value = 4733.1665
$ws.Range("M62").Value = -4109.1665
$ws.Range("H65").Value = 6151.375
$ws.Range("I65").Value = 4733.1665
$ws.Range("K65").Value = 23665.8325
$ws.Range("M65").Value = -20545.8325
$ws.Range("H86").Value = 7100.4443
$ws.Range("I86").Value = 7162.5
$ws.Range("K86").Value = 7162.5
$ws.Range("M86").Value = -6039.5
$ws.Range("H89").Value = 7100.4443
$ws.Range("I89").Value = 7162.5
$ws.Range("K89").Value = 35812.5
$ws.Range("M89").Value = -30196.5
$ws.Range("H98").Value = 2062.875
$ws.Range("I98").Value = 1214.7142
$ws.Range("J98").Value = 8000
$ws.Range("K98").Value = 1214.7142
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = 283.2858000000001
$ws.Range("N98").Value = -10996
$ws.Range("H106").Value = 8563.305
$ws.Range("J106").Value = 15541.3
$ws.Range("L106").Value = 15541.3
$ws.Range("N106").Value = -16803.3
$ws.Range("H116").Value = 12673
$ws.Range("I116").Value = 9649.5
$ws.Range("J116").Value = 16301.2
$ws.Range("K116").Value = 9649.5
$ws.Range("L116").Value = 16301.2
$ws.Range("M116").Value = -6207.5
$ws.Range("N116").Value = -23185.2
$ws.Range("H122").Value = 2062.875
$ws.Range("I122").Value = 1214.7142
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 3644.1426
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -1194.1426
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 1649.3784
$ws.Range("I132").Value = 1425.2122
$ws.Range("K132").Value = 4275.6366
$ws.Range("M132").Value = -1745.6366

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2734.8914
$ws.Range("I32").Value = 1851
$ws.Range("K32").Value = 1851
$ws.Range("M32").Value = -1564
$ws.Range("H61").Value = 3693.318
$ws.Range("I61").Value = 2929.158
$ws.Range("K61").Value = 2929.158
$ws.Range("M61").Value = -2717.158
$ws.Range("H74").Value = 15152986
$ws.Range("I74").Value = 22223566
$ws.Range("K74").Value = 22223566
$ws.Range("M74").Value = -22222692
$ws.Range("H77").Value = 15152986
$ws.Range("I77").Value = 22223566
$ws.Range("K77").Value = 111117830
$ws.Range("M77").Value = -111113462
$ws.Range("H97").Value = 861.5
$ws.Range("I97").Value = 480.5
$ws.Range("K97").Value = 480.5
$ws.Range("M97").Value = 15.5
$ws.Range("H122").Value = 3369.7778
$ws.Range("I122").Value = 1645.5834
$ws.Range("J122").Value = 6818.1665
$ws.Range("K122").Value = 4936.7502
$ws.Range("L122").Value = 20454.4995
$ws.Range("M122").Value = -2486.7502
$ws.Range("N122").Value = -25354.4995
$ws.Range("H129").Value = 174997
$ws.Range("J129").Value = 174997
$ws.Range("L129").Value = 174997
$ws.Range("N129").Value = -184997
$ws.Range("H136").Value = 3693.318
$ws.Range("I136").Value = 2929.158
$ws.Range("K136").Value = 8787.474
$ws.Range("M136").Value = -6237.474

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3875.5833
$ws.Range("I20").Value = 2560.8
$ws.Range("J20").Value = 4814.7144
$ws.Range("K20").Value = 2560.8
$ws.Range("L20").Value = 4814.7144
$ws.Range("M20").Value = -2313.8
$ws.Range("N20").Value = -5308.7144
$ws.Range("H86").Value = 1656.7142
$ws.Range("I86").Value = 1349.5
$ws.Range("K86").Value = 1349.5
$ws.Range("M86").Value = -226.5
$ws.Range("H89").Value = 1656.7142
$ws.Range("I89").Value = 1349.5
$ws.Range("K89").Value = 6747.5
$ws.Range("M89").Value = -1131.5
$ws.Range("H99").Value = 1676.3334
$ws.Range("I99").Value = 1170.7142
$ws.Range("K99").Value = 1170.7142
$ws.Range("M99").Value = 327.2858000000001
$ws.Range("H107").Value = 787.5
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 800
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1120
$ws.Range("N107").Value = -4590

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37432.812
$ws.Range("J31").Value = 71114.56
$ws.Range("L31").Value = 71114.56
$ws.Range("N31").Value = -71704.56
$ws.Range("H34").Value = 37432.812
$ws.Range("J34").Value = 71114.56
$ws.Range("L34").Value = 71114.56
$ws.Range("N34").Value = -71518.56
$ws.Range("H58").Value = 4246.892
$ws.Range("J58").Value = 6905.222
$ws.Range("L58").Value = 6905.222
$ws.Range("N58").Value = -7311.222
$ws.Range("H107").Value = 1416.2354
$ws.Range("I107").Value = 1392.6428
$ws.Range("K107").Value = 1392.6428
$ws.Range("M107").Value = 527.3571999999999
$ws.Range("H136").Value = 4246.892
$ws.Range("J136").Value = 6905.222
$ws.Range("L136").Value = 20715.666
$ws.Range("N136").Value = -25815.666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 7072.5
$ws.Range("I68").Value = 4650
$ws.Range("K68").Value = 13950
$ws.Range("M68").Value = -13139
$ws.Range("H71").Value = 7072.5
$ws.Range("I71").Value = 4650
$ws.Range("K71").Value = 41850
$ws.Range("M71").Value = -37794
$ws.Range("H132").Value = 4553.6
$ws.Range("I132").Value = 3421.2856
$ws.Range("J132").Value = 5544.375
$ws.Range("K132").Value = 30791.5704
$ws.Range("L132").Value = 49899.375
$ws.Range("M132").Value = -28261.5704
$ws.Range("N132").Value = -54959.375
$ws.Range("H134").Value = 6103.3687
$ws.Range("I134").Value = 2980.9167
$ws.Range("K134").Value = 8942.750100000001
$ws.Range("M134").Value = -3872.750100000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1112.5385
$ws.Range("I2").Value = 203.33333
$ws.Range("K2").Value = 203.33333
$ws.Range("M2").Value = -90.33332999999999
$ws.Range("H97").Value = 1924.1765
$ws.Range("I97").Value = 1515.8334
$ws.Range("K97").Value = 1515.8334
$ws.Range("M97").Value = -1019.8334
$ws.Range("H102").Value = 3097.75
$ws.Range("I102").Value = 1640.5555
$ws.Range("K102").Value = 1640.5555
$ws.Range("M102").Value = -18.55549999999994
$ws.Range("H126").Value = 3300.7
$ws.Range("I126").Value = 1656.8572
$ws.Range("K126").Value = 4970.571599999999
$ws.Range("M126").Value = -2500.571599999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3092.5
$ws.Range("I16").Value = 2320.1428
$ws.Range("K16").Value = 2320.1428
$ws.Range("M16").Value = -2150.1428
$ws.Range("H22").Value = 4113.1333
$ws.Range("I22").Value = 2140
$ws.Range("J22").Value = 5099.7
$ws.Range("K22").Value = 2140
$ws.Range("L22").Value = 5099.7
$ws.Range("M22").Value = -1845
$ws.Range("N22").Value = -5689.7
$ws.Range("H27").Value = 4113.1333
$ws.Range("I27").Value = 2140
$ws.Range("J27").Value = 5099.7
$ws.Range("K27").Value = 2140
$ws.Range("L27").Value = 5099.7
$ws.Range("M27").Value = -2033
$ws.Range("N27").Value = -5313.7
$ws.Range("H46").Value = 3683.5
$ws.Range("I46").Value = 3400.3333
$ws.Range("J46").Value = 3966.6667
$ws.Range("K46").Value = 3400.3333
$ws.Range("L46").Value = 3966.6667
$ws.Range("M46").Value = -3212.3333
$ws.Range("N46").Value = -4342.6667
$ws.Range("H136").Value = 4203.914
$ws.Range("I136").Value = 2010.96
$ws.Range("K136").Value = 6032.88
$ws.Range("M136").Value = -3482.88

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1410.5555
$ws.Range("I100").Value = 598.1667
$ws.Range("K100").Value = 1196.3334
$ws.Range("M100").Value = -655.3334
$ws.Range("H107").Value = 999.6667
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 12343.333
$ws.Range("I122").Value = 3361.6667
$ws.Range("K122").Value = 10085.0001
$ws.Range("M122").Value = -7635.000100000001
$ws.Range("H132").Value = 2700.0205
$ws.Range("I132").Value = 2751
$ws.Range("J132").Value = 2438.75
$ws.Range("K132").Value = 8253
$ws.Range("L132").Value = 7316.25
$ws.Range("M132").Value = -5723
$ws.Range("N132").Value = -12376.25
$ws.Range("H136").Value = 7557.4287
$ws.Range("J136").Value = 10100.75
$ws.Range("L136").Value = 30302.25
$ws.Range("N136").Value = -35402.25
